$wb = $excel.ActiveWorkbook

# Sheet 1 contains the metadata table (Property / Value)
$ws1 = $wb.Worksheets.Item(1)

# Row 4 is "Name" - set its Value cell (B4) which was previously empty
$ws1.Range("B4").Value = "CnamamelisecteurconventionnementVs"

# Row 8 is "Date" - update its Value cell (B8) to the new timestamp
$ws1.Range("B8").Value = "2025-07-18T06:40:38+00:00"
